# Re-style the three tables that used the custom "Table_0" style
# ({0DEEB2D3-7035-4EE4-A8E2-CD3FC609551B}) so they use the built-in
# table style {7FDC08DD-1170-4111-99A6-D6BBD1FB67CA} instead.
#
# The tables live on slides 14, 15 and 16, each as the first shape
# on the slide (a graphicFrame hosting the table).

$p = $ppt.ActivePresentation

$oldStyleId = "{0DEEB2D3-7035-4EE4-A8E2-CD3FC609551B}"
$newStyleId = "{7FDC08DD-1170-4111-99A6-D6BBD1FB67CA}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
